# Added ability to strip out * from import headers
# Update the import header row so that header labels are suffixed with " *"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name *"
$ws.Range("C1").Value = "Category *"
$ws.Range("E1").Value = "Fund *"

$ws.Range("E1:E7").Style = "Normal 2"

$ws.Range("C2").Select()
